$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CALCULO")
$ws.Range("A1").Value = "test"
